$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text format on price cells whose new values would otherwise be
# auto-converted to numbers by Excel (single-dot decimal-looking strings),
# so the literal text (incl. trailing zeros) is preserved like the original inlineStr cells.
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D48").NumberFormat = "@"

# Apply the updated cryptos list values
$ws.Range("D2").Value = "36.677.68"
$ws.Range("E2").Value = "  -0.78%  "
$ws.Range("D3").Value = "2.058.01"
$ws.Range("E3").Value = "  +0.94%  "
$ws.Range("E4").Value = "  +0.02%  "
$ws.Range("D5").Value = "243.60"
$ws.Range("E5").Value = "  -0.71%  "
$ws.Range("D6").Value = "0.666"
$ws.Range("E6").Value = "  +1.78%  "
$ws.Range("D8").Value = "54.56"
$ws.Range("E8").Value = "  -6.66%  "
$ws.Range("D9").Value = "59.23"
$ws.Range("E9").Value = "  +0.75%  "
$ws.Range("E10").Value = "  -2.74%  "
$ws.Range("E12").Value = "  -3.04%  "
$ws.Range("D13").Value = "0.934"
$ws.Range("E13").Value = "  +5.35%  "
$ws.Range("D14").Value = "14.77"
$ws.Range("D15").Value = "2.358.54"
$ws.Range("E15").Value = "  +0.96%  "
$ws.Range("E16").Value = "  -3.17%  "
$ws.Range("D17").Value = "2.063.81"
$ws.Range("E17").Value = "  +4.03%  "
$ws.Range("D18").Value = "36.582.69"
$ws.Range("E18").Value = "  -0.98%  "
$ws.Range("D19").Value = "17.09"
$ws.Range("E19").Value = "  -6.44%  "
$ws.Range("D20").Value = "72.03"
$ws.Range("E20").Value = "  -1.96%  "
$ws.Range("D22").Value = "238.10"
$ws.Range("E22").Value = "  +1.26%  "
$ws.Range("E23").Value = "  -1.74%  "
$ws.Range("E24").Value = "  -0.06%  "
$ws.Range("E25").Value = "  -3.16%  "
$ws.Range("E26").Value = "  -0.22%  "
$ws.Range("D27").Value = "9.35"
$ws.Range("E27").Value = "  -2.11%  "
$ws.Range("D28").Value = "164.61"
$ws.Range("E28").Value = "  -1.78%  "
$ws.Range("D29").Value = "20.15"
$ws.Range("E29").Value = "  +1.33%  "
$ws.Range("E30").Value = "  -1.33%  "
$ws.Range("B31").Value = "ImmutableX"
$ws.Range("C31").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D31").Value = "1.20"
$ws.Range("E31").Value = "  +7.83%  "
$ws.Range("B32").Value = "Filecoin"
$ws.Range("C32").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D32").Value = "5.09"
$ws.Range("E32").Value = "  -7.56%  "
$ws.Range("D33").Value = "4.52"
$ws.Range("E33").Value = "  -4.81%  "
$ws.Range("E34").Value = "  -2.21%  "
$ws.Range("E35").Value = "  +0.03%  "
$ws.Range("E36").Value = "  -0.59%  "
$ws.Range("D37").Value = "2.22"
$ws.Range("E37").Value = "  -0.87%  "
$ws.Range("D38").Value = "0.0836"
$ws.Range("E38").Value = "  -2.28%  "
$ws.Range("E39").Value = "  -3.47%  "
$ws.Range("D40").Value = "4.96"
$ws.Range("E40").Value = "  -4.93%  "
$ws.Range("D41").Value = "2.89"
$ws.Range("E41").Value = "  -6.64%  "
$ws.Range("E42").Value = "  -2.50%  "
$ws.Range("E43").Value = "  -2.75%  "
$ws.Range("D44").Value = "94.40"
$ws.Range("E44").Value = "  -2.52%  "
$ws.Range("D45").Value = "0.0910"
$ws.Range("E45").Value = "  -3.74%  "
$ws.Range("D46").Value = "1.411.53"
$ws.Range("E46").Value = "  +9.24%  "
$ws.Range("D47").Value = "7.67"
$ws.Range("E47").Value = "  +15.15%  "
$ws.Range("D48").Value = "16.04"
$ws.Range("E48").Value = "  -4.70%  "
$ws.Range("E49").Value = "  +2.15%  "
$ws.Range("E50").Value = "  -3.84%  "
$ws.Range("D51").Value = "2.246.77"
$ws.Range("E51").Value = "  +1.09%  "
